$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend style (bold, centered, bordered) used in column A for rows 10-16 down through new rows 17-19
$ws.Range("A10").Copy()
$ws.Range("A11:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10: Gaussian-Quadrature (moved up from old row 16, same averaged-intensity data)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.175415591871771
$ws.Range("D10").Value = 1.160408762032967
$ws.Range("E10").Value = 0.9646712100274477
$ws.Range("F10").Value = 0.9831432012714794
$ws.Range("G10").Value = 1.175415591871771
$ws.Range("H10").Value = 1.160408762032967
$ws.Range("I10").Value = 0.9601667918494853
$ws.Range("J10").Value = 0.9213191143260228
$ws.Range("K10").Value = 1.038141564909271
$ws.Range("L10").Value = 0.9832396425783597
$ws.Range("M10").Value = 1.174892722463677
$ws.Range("N10").Value = 1.062539986030207
$ws.Range("O10").Value = 1.070909691300916
$ws.Range("P10").Value = 1.02331323485835

# Row 11: NEW: Spiral-90deg-10rot-5space
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.8792296606471938
$ws.Range("D11").Value = 0.8901937530894197
$ws.Range("E11").Value = 1.223956596199306
$ws.Range("F11").Value = 0.991572493718043
$ws.Range("G11").Value = 0.8792296606471938
$ws.Range("H11").Value = 0.8901937530894197
$ws.Range("I11").Value = 1.005462893250961
$ws.Range("J11").Value = 1.056593959752186
$ws.Range("K11").Value = 0.949508368190786
$ws.Range("L11").Value = 0.9169978765271285
$ws.Range("M11").Value = 0.8786506806244286
$ws.Range("N11").Value = 1.057075174644363
$ws.Range("O11").Value = 0.9962381259134907
$ws.Range("P11").Value = 0.9891894501718781

# Row 12: NEW: Spiral-90deg-15rot-5space
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.8775495200646271
$ws.Range("D12").Value = 0.8899116419768084
$ws.Range("E12").Value = 1.224593306428732
$ws.Range("F12").Value = 0.9920955531385723
$ws.Range("G12").Value = 0.8775495200646271
$ws.Range("H12").Value = 0.8899116419768084
$ws.Range("I12").Value = 1.005150463758199
$ws.Range("J12").Value = 1.057065762041839
$ws.Range("K12").Value = 0.9492849771673302
$ws.Range("L12").Value = 0.9173206177966033
$ws.Range("M12").Value = 0.8769685892509357
$ws.Range("N12").Value = 1.05725247420277
$ws.Range("O12").Value = 0.9960375054021851
$ws.Range("P12").Value = 0.9891214802965891

# Row 13: NEW: Spiral-90deg-10rot-3space
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.8788903103605864
$ws.Range("D13").Value = 0.8900734918225819
$ws.Range("E13").Value = 1.223959563745626
$ws.Range("F13").Value = 0.9917704174351063
$ws.Range("G13").Value = 0.8788903103605864
$ws.Range("H13").Value = 0.8900734918225819
$ws.Range("I13").Value = 1.005334164345989
$ws.Range("J13").Value = 1.056715228505023
$ws.Range("K13").Value = 0.9495046666451685
$ws.Range("L13").Value = 0.9170780259792974
$ws.Range("M13").Value = 0.8783098151014727
$ws.Range("N13").Value = 1.057016527784104
$ws.Range("O13").Value = 0.9961734458409751
$ws.Range("P13").Value = 0.9891657336049222

# Row 14: NoRotation-tilt60deg (was row 10)
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.5651960000000001
$ws.Range("D14").Value = 0.7293319999999983
$ws.Range("E14").Value = 1.792452000000003
$ws.Range("F14").Value = 0.9803519999999986
$ws.Range("G14").Value = 0.5651960000000001
$ws.Range("H14").Value = 0.7293319999999983
$ws.Range("I14").Value = 1.047347999999999
$ws.Range("J14").Value = 1.157292000000001
$ws.Range("K14").Value = 0.8337279999999985
$ws.Range("L14").Value = 0.8057599999999999
$ws.Range("M14").Value = 0.5647799999999997
$ws.Range("N14").Value = 1.260892000000001
$ws.Range("O14").Value = 1.016833
$ws.Range("P14").Value = 0.9889324999999998

# Row 15: Rotation-NoTilt (was row 11)
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.45
$ws.Range("D15").Value = 0.61
$ws.Range("E15").Value = 2.177237500000003
$ws.Range("F15").Value = 0.93
$ws.Range("G15").Value = 0.45
$ws.Range("H15").Value = 0.61
$ws.Range("I15").Value = 1.101662499999999
$ws.Range("J15").Value = 1.18
$ws.Range("K15").Value = 0.77
$ws.Range("L15").Value = 0.6899999999999999
$ws.Range("M15").Value = 0.45
$ws.Range("N15").Value = 1.393618750000002
$ws.Range("O15").Value = 1.041809375000001
$ws.Range("P15").Value = 0.9886125000000004

# Row 16: Rotation-60detTilt (was row 12)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.6799373139968016
$ws.Range("D16").Value = 0.7671106207744022
$ws.Range("E16").Value = 1.660935064780803
$ws.Range("F16").Value = 0.953186363904
$ws.Range("G16").Value = 0.6799373139968016
$ws.Range("H16").Value = 0.7671106207744022
$ws.Range("I16").Value = 1.061125174374399
$ws.Range("J16").Value = 1.100664622489596
$ws.Range("K16").Value = 0.8649486254080015
$ws.Range("L16").Value = 0.8200762056703989
$ws.Range("M16").Value = 0.6799242782720016
$ws.Range("N16").Value = 1.214022842777603
$ws.Range("O16").Value = 1.015292340864002
$ws.Range("P16").Value = 0.9884979989248003

# Row 17: HexGrid-90degTilt5degRes (was row 13)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9925198466865341
$ws.Range("D17").Value = 0.9923480365276607
$ws.Range("E17").Value = 0.9912463130490972
$ws.Range("F17").Value = 0.9909893856402903
$ws.Range("G17").Value = 0.9925198466865341
$ws.Range("H17").Value = 0.9923480365276607
$ws.Range("I17").Value = 0.990870661866622
$ws.Range("J17").Value = 0.9908381835054233
$ws.Range("K17").Value = 0.9905033285288765
$ws.Range("L17").Value = 0.9904279454116482
$ws.Range("M17").Value = 0.9925003681086734
$ws.Range("N17").Value = 0.9917971747883789
$ws.Range("O17").Value = 0.9917758954758955
$ws.Range("P17").Value = 0.991217962652019

# Row 18: HexGrid-90degTilt22p5degRes (was row 14)
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9568823335224449
$ws.Range("D18").Value = 0.9731016782944861
$ws.Range("E18").Value = 1.030981562966724
$ws.Range("F18").Value = 0.9937621446945265
$ws.Range("G18").Value = 0.9568823335224449
$ws.Range("H18").Value = 0.9731016782944861
$ws.Range("I18").Value = 0.9788556014680605
$ws.Range("J18").Value = 0.9905500323813348
$ws.Range("K18").Value = 0.9866794384209309
$ws.Range("L18").Value = 0.9838975430333876
$ws.Range("M18").Value = 0.9568823335224449
$ws.Range("N18").Value = 1.002041620630605
$ws.Range("O18").Value = 0.9886819298695454
$ws.Range("P18").Value = 0.9868387918477369

# Row 19: HexGrid-60degTilt5degRes (was row 15)
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.999519819632665
$ws.Range("D19").Value = 1.017562233915762
$ws.Range("E19").Value = 0.9828356523002032
$ws.Range("F19").Value = 0.9922549711631075
$ws.Range("G19").Value = 0.999519819632665
$ws.Range("H19").Value = 1.017562233915762
$ws.Range("I19").Value = 0.9834194394324258
$ws.Range("J19").Value = 0.9833850312650245
$ws.Range("K19").Value = 0.988071276157061
$ws.Range("L19").Value = 1.002057048397828
$ws.Range("M19").Value = 0.9994806574392824
$ws.Range("N19").Value = 1.000198943107983
$ws.Range("O19").Value = 0.9980431692529343
$ws.Range("P19").Value = 0.9936381840330096

Write-Host "Done"
